$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - add new header "Id" in column F
$ws.Range("F1").Value = "Id"

# Row 2 - overwrite text fields with "asd"; E2 is left untouched (unchanged in target);
# F2 is a new numeric cell
$ws.Range("A2").Value = "asd"
$ws.Range("B2").Value = "asd"
$ws.Range("C2").Value = "asd"
$ws.Range("D2").Value = "asd"
$ws.Range("F2").Value = 1

# Row 3 - new row with text fields "qwe" and numeric fields
$ws.Range("A3").Value = "qwe"
$ws.Range("B3").Value = "qwe"
$ws.Range("C3").Value = "qwe"
$ws.Range("D3").Value = "qwe"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
